# Applies the LOT2012 content edits:
#  - inserts a new row at row 13 (shifts old rows 13-23 down to 14-24)
#  - fills in previously-missing "B"/"C" text for several rows
#  - fixes the newly inserted row's formatting to match column B/C styles

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row at position 13; everything below shifts down one row.
$ws.Rows.Item(13).Insert()

# The insert carried column A's (bold) formatting across the whole new row;
# clear it away, then copy B/C's normal + red styles down from the row below
# (row 14, which is the old row 13) onto the new row 13.
$ws.Range("A13").Clear()
$ws.Range("B14:C14").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false


$v = 'Introdução à Engenharia Ambiental para estudantes de Engenharia Bioquímica. Dar conhecimentos aos alunos de noções básicas sobre ecologia e impacto das atividades da engenharia sobre o meio ambiente Estudo da preservação do meio ambiente e tratamentos aplicados aos poluentes locais e globais. Conceitos legais e institucionais para o desenvolvimento sustentável.'
$ws.Range("B10").Value2 = $v
$ws.Range("C10").Value2 = $v

$v = '1720367 - Teresa Cristina Brazil de Paiva'
$ws.Range("B13").Value2 = $v
$ws.Range("C13").Value2 = $v

$v = '1-Fundamentos; 2- Poluição Ambiental; 3- Desenvolvimento Sustentável'
$ws.Range("B14").Value2 = $v
$ws.Range("C14").Value2 = $v

$v = '1,FUNDAMENTOS: A crise ambiental e as leis da física. Fluxo de Energia nos ecossistemas, cadeias alimentares, sucessão ecológica e ciclos biogeoquímicos. Dinâmica das populações. Base para o desenvilmento sustentável. 2, POLUIÇÃO AMBIENTAL, CONTROLE E TRATAMENTO: O conceito de poluição e seu controle. O Meio Aquático: usos e requisitos de qualidades das águas parâmetros característicos da água. Poluição: fontes e poluição biodegradação, poluentes tóxicos e metais pesados, comportamento dos poluentes no meio aquático. Poluição em lagos: estratificação térmica e eutrofização , monitoramento da poluição da água, poluição difusa urbana e rural, Tratamento da água e de esgotos. O Meio Terrestre: origem, composição e formação dos solos, erosão e seu controle. Poluição do solo rural: fertilizantes, defensivos agrícolas, formas alternativas de controle de pragas do solo urbano, Formas de disposição e tratamento do lixo urbarno: compostagem, incineração e aterro sanitário. Resíduos perigosos: fontes, efeitos sobre a saúde e disposição e tratamento. O Meio Atmosférico: poluição global efeito estufa e camada de ozônio. Poluição local e regional: smog industrial e fotoquímico, efeitos da poluição do ar. Meteorologia e dispersão de poluentes: Processo de dispersão de plumas, controle da poluição do ar nas grandes cidades brasileiras. Poluição sonora. 3. DESENVOLVIMENTO SUSTENTÁVEL: A crise energética, fontes alternativas de energia. , análise econômica, relação benefício-custo, instrumentos de planejamento e gestão: alocação de custos, cobrança pelo uso dos recursos naturais, principio poluidor-pagador, outorga de uso de recursos naturais, avaliação de impactos ambientais: descrição geral, indicadores de impacto, métodos quantitativos, aspectos legais e institucionais do controle ambiental.'
$ws.Range("B16").Value2 = $v
$ws.Range("C16").Value2 = $v

$v = 'Os alunos serão avaliados por meio de duas provas (P1 e P2) e complementada por meio de trabalhos, seminários e/ou relatórios (C).'
$ws.Range("B19").Value2 = $v
$ws.Range("C19").Value2 = $v

$v = 'A nota final (NF) será calculada atribuindo-se peso um para a primeira avaliação (P1 = 7 pontos e C = 3 pontos) e peso dois para a segunda avaliação (P2 = 10 pontos).A média ponderada das notas corresponderá à média do período letivo, ou seja: Média do período letivo normal = ((P1 + C) + P2.2)/3.Serão aprovados os alunos que obtiverem média igual ou maior que 5,0 e 70% de frequência no curso.'
$ws.Range("B20").Value2 = $v
$ws.Range("C20").Value2 = $v

$v = 'Aos alunos que obtiverem média igual ou maior que 3,0 e menor que 5,0 será oferecido um programa de recuperação que será avaliado por uma prova final. Nesse caso, a média final do aluno será: Média final = (média do período letivo normal + nota prova final)/2.Serão aprovados os alunos que obtiverem média final igual ou maior que 5,0.'
$ws.Range("B21").Value2 = $v
$ws.Range("C21").Value2 = $v

$v = '1,Braga, B.P.F., M.T.,Conejo, J.G., Porto, M.F., Veras M.S., Nucci, N., Juliano, N. e Eiger, S. - Introdução à Engenharia Ambiental, Makron Books, São Paulo, 1998, 2. Sperling, M.V. - Princípios do Tratamento Biológico de Águas Residuárias. Desa-UFMG, Minas Gerais, 1996.'
$ws.Range("B22").Value2 = $v
$ws.Range("C22").Value2 = $v
